$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Extend the year table with a new 2023 column (T), copying the formatting
# of the existing 2022 column (S) so the new header/value cells look the
# same as their neighbours.
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial($xlPasteFormats)
$ws.Range("T4").Value = 2023

$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial($xlPasteFormats)
$ws.Range("T5").Value = 75.099999999999994

$excel.CutCopyMode = 0

# Adjust column widths: A:C got slightly narrower, D:T got an explicit width
# to fit the now-20-column-wide table.
$ws.Range("A1:C1").ColumnWidth = 35.5703125
$ws.Range("D1:T1").ColumnWidth = 8.85546875

# Put the active selection back on A1 so the saved view no longer points at
# the stale P8 cell.
$ws.Range("A1").Select()
